$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the existing wave (row 3): rename from the old "OCI Trial" sample
# to the new "First Flow" Linux leg, keeping its IP / OS / username / sync type.
$ws.Range("A3").Value = "First Flow"
$ws.Range("E3").Value = "psp-MyLinFirstFlow-src1"
$ws.Range("J3").Value = "psp-MyLinFirstFlow-tgt1"

# --- Add the new Windows leg of the first automation flow as row 4.
$ws.Range("A4").Value = "First Flow"
$ws.Range("B4").Value = $true
$ws.Range("C4").Value = "Autoprovision"
$ws.Range("D4").Value = "172.29.30.196"
$ws.Range("E4").Value = "psp-MyWinFirstFlow-src1"
$ws.Range("F4").Value = "Windows"
$ws.Range("G4").Value = 22
$ws.Range("H4").Value = "SYSTEM"
$ws.Range("I4").Value = "Direct Sync"
$ws.Range("J4").Value = "psp-MyWinFirstFlow-tgt1"

# --- Resize columns whose best-fit width changed because of the new/longer text
# (mirrors Excel's automatic "best fit" recompute after the new values were entered).
$ws.Columns.Item(1).ColumnWidth = 16.109375
$ws.Columns.Item(3).ColumnWidth = 13.44140625
$ws.Columns.Item(5).ColumnWidth = 23.21875
$ws.Columns.Item(10).ColumnWidth = 23.109375
$ws.Columns.Item(11).ColumnWidth = 18.77734375

# --- Move the selection the way it was left after the edit.
$ws.Range("J4").Select() | Out-Null
